# Fruta / hortaliza, semanal
# Inserts a new week's worth of "Brócoli" price records (2 rows) ahead of the
# existing data block, and appends another week's worth (2 rows) at the end,
# mirroring the recurring per-date, Primera/Segunda row-pair pattern already
# present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two fresh rows at the top of the date block (old row 81),
#        pushing the existing rows 81:128 down to 83:130. ---
$ws.Range("A81:R82").EntireRow.Insert()

# Row 81 - Calidad "Primera"
$ws.Cells.Item(81, 1).Value = 11
$ws.Cells.Item(81, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(81, 3).Value = "Bíobío"
$ws.Cells.Item(81, 4).Value = 44435
$ws.Cells.Item(81, 5).Value = 8
$ws.Cells.Item(81, 6).Value = 100112023
$ws.Cells.Item(81, 7).Value = "Brócoli"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 3000
$ws.Cells.Item(81, 11).Value = 600
$ws.Cells.Item(81, 12).Value = 800
$ws.Cells.Item(81, 13).Value = 683
$ws.Cells.Item(81, 14).Value = "$/unidad"
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 683
$ws.Cells.Item(81, 17).Value = 1
$ws.Cells.Item(81, 18).Value = "Hortaliza"

# Row 82 - Calidad "Segunda"
$ws.Cells.Item(82, 1).Value = 11
$ws.Cells.Item(82, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(82, 3).Value = "Bíobío"
$ws.Cells.Item(82, 4).Value = 44435
$ws.Cells.Item(82, 5).Value = 8
$ws.Cells.Item(82, 6).Value = 100112023
$ws.Cells.Item(82, 7).Value = "Brócoli"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Segunda"
$ws.Cells.Item(82, 10).Value = 1500
$ws.Cells.Item(82, 11).Value = 500
$ws.Cells.Item(82, 12).Value = 600
$ws.Cells.Item(82, 13).Value = 533
$ws.Cells.Item(82, 14).Value = "$/unidad"
$ws.Cells.Item(82, 15).Value = "Región Metropolitana"
$ws.Cells.Item(82, 16).Value = 533
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = "Hortaliza"

# --- 2. Append two more rows (now 131:132) with another week's records. ---

# Row 131 - Calidad "Primera"
$ws.Cells.Item(131, 1).Value = 11
$ws.Cells.Item(131, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(131, 3).Value = "Bíobío"
$ws.Cells.Item(131, 4).Value = 44432
$ws.Cells.Item(131, 5).Value = 8
$ws.Cells.Item(131, 6).Value = 100112023
$ws.Cells.Item(131, 7).Value = "Brócoli"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 2000
$ws.Cells.Item(131, 11).Value = 600
$ws.Cells.Item(131, 12).Value = 700
$ws.Cells.Item(131, 13).Value = 650
$ws.Cells.Item(131, 14).Value = "$/unidad"
$ws.Cells.Item(131, 15).Value = "Región Metropolitana"
$ws.Cells.Item(131, 16).Value = 650
$ws.Cells.Item(131, 17).Value = 1
$ws.Cells.Item(131, 18).Value = "Hortaliza"

# Row 132 - Calidad "Segunda"
$ws.Cells.Item(132, 1).Value = 11
$ws.Cells.Item(132, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(132, 3).Value = "Bíobío"
$ws.Cells.Item(132, 4).Value = 44432
$ws.Cells.Item(132, 5).Value = 8
$ws.Cells.Item(132, 6).Value = 100112023
$ws.Cells.Item(132, 7).Value = "Brócoli"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Segunda"
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 500
$ws.Cells.Item(132, 12).Value = 500
$ws.Cells.Item(132, 13).Value = 500
$ws.Cells.Item(132, 14).Value = "$/unidad"
$ws.Cells.Item(132, 15).Value = "Región Metropolitana"
$ws.Cells.Item(132, 16).Value = 500
$ws.Cells.Item(132, 17).Value = 1
$ws.Cells.Item(132, 18).Value = "Hortaliza"

# Give the new date cells (column D) the same date style used throughout
# the rest of the column (style index 2 on the existing sheet).
$ws.Range("D81:D82").NumberFormat = $ws.Range("D80").NumberFormat
$ws.Range("D131:D132").NumberFormat = $ws.Range("D130").NumberFormat
